# Refresh cryptos list: update Price (D) and Volume(1h) (E) columns to match
# the latest scrape snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.894.81'
$ws.Range("E2").Value = '  +3.03%  '
$ws.Range("D3").Value = '1.785.18'
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("E4").Value = '  +0.36%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '225.29'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.559'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.63%  '
$ws.Range("E8").Value = '  -4.81%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.36'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +3.34%  '
$ws.Range("E10").Value = '  +0.25%  '
$ws.Range("E11").Value = '  -0.29%  '
$ws.Range("E12").Value = '  +0.61%  '
$ws.Range("D13").Value = '2.042.85'
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '1.780.78'
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '10.45'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("E16").Value = '  -1.26%  '
$ws.Range("D17").Value = '33.920.28'
$ws.Range("E17").Value = '  +3.35%  '
$ws.Range("E18").Value = '  -3.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.11'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.49%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '252.04'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.94%  '
$ws.Range("D21").Value = '0.0₃0740'
$ws.Range("E21").Value = '  -0.77%  '
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.33'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.86%  '
$ws.Range("E24").Value = '  -3.16%  '
$ws.Range("E25").Value = '  -1.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.14'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.48'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -0.48%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.99'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.06%  '
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("E32").Value = '  -0.34%  '
$ws.Range("E33").Value = '  +1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.60'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +2.74%  '
$ws.Range("E35").Value = '  +3.21%  '
$ws.Range("D36").Value = '1.503.06'
$ws.Range("E36").Value = '  -3.68%  '
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.632'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.18%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '83.46'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.04%  '
$ws.Range("E40").Value = '  -0.22%  '
$ws.Range("E41").Value = '  +1.48%  '
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("E43").Value = '  +2.42%  '
$ws.Range("E44").Value = '  +0.97%  '
$ws.Range("E45").Value = '  -2.03%  '
$ws.Range("E46").Value = '  +3.03%  '
$ws.Range("D47").Value = '1.936.13'
$ws.Range("E47").Value = '  +1.34%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.72'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.71%  '
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("E50").Value = '  +8.27%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '51.45'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -4.98%  '
